# Weekly update to the "Poroto verde" (green bean) price sheet:
# two new rows of price data were inserted right after the row for
# 2020-12-02 (row 153) / before the old row 154, pushing everything
# below down by two rows (old 154..169 -> new 156..171).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 154/155 (rows 154..169 shift to 156..171)
$ws.Rows.Item(154).Resize(2).EntireRow.Insert()

# --- New row 154 -----------------------------------------------------
$ws.Cells.Item(154, 1).Value = 8
$ws.Cells.Item(154, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(154, 3).Value = "Coquimbo"
$ws.Cells.Item(154, 4).Value = 44578
$ws.Cells.Item(154, 5).Value = 4
$ws.Cells.Item(154, 6).Value = 100112031
$ws.Cells.Item(154, 7).Value = "Poroto verde"
$ws.Cells.Item(154, 8).Value = "Magnum"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 400
$ws.Cells.Item(154, 11).Value = 25000
$ws.Cells.Item(154, 12).Value = 26000
$ws.Cells.Item(154, 13).Value = 25500
$ws.Cells.Item(154, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(154, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(154, 16).Value = 1020
$ws.Cells.Item(154, 17).Value = 25
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# --- New row 155 -----------------------------------------------------
$ws.Cells.Item(155, 1).Value = 8
$ws.Cells.Item(155, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44578
$ws.Cells.Item(155, 5).Value = 4
$ws.Cells.Item(155, 6).Value = 100112031
$ws.Cells.Item(155, 7).Value = "Poroto verde"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 400
$ws.Cells.Item(155, 11).Value = 27000
$ws.Cells.Item(155, 12).Value = 28000
$ws.Cells.Item(155, 13).Value = 27500
$ws.Cells.Item(155, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(155, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(155, 16).Value = 1100
$ws.Cells.Item(155, 17).Value = 25
$ws.Cells.Item(155, 18).Value = "Hortaliza"

Write-Output "Inserted rows 154-155 and populated new data"
